# Applies the LOM3049.docx content updates described by the diff.
$d = $word.ActiveDocument

# 1) Update activation date.
$d.Content.Find.Execute(
    "Ativação: 01/01/2021", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2024", 2)

# 2) Replace the "Objetivos" paragraph text.
$d.Content.Find.Execute(
    "Abordar os princípios básicos da termodinâmica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e sólido sobre estes princípios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodinâmica é aplicada na prática de engenharia. Enfatizar a compreensão da termodinâmica baseada na Física e em argumentos físicos, buscando incentivar o entendimento mais profundo da termodinâmica.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Esta disciplina faz parte da formação do engenheiro de materiais, contribuindo para gerar competências gerais e específicas.Abordar os princípios básicos da Termodinâmica dentro do contexto de máquinas térmicas.Incentivar os alunos a identificar como a termodinâmica está relacionada com as principais atividades humanas, com ênfase na geração de potência e refrigeração.Relacionar esta disciplina com outras da grade do curso, como: Física, Recursos Naturais, Tecnologias Limpas para Geração de Energia, Termodinâmica de Materiais, Seleção de Materiais, Fenômenos de Transporte p/ EM, dentre outras. Desenvolver nos alunos a prática da busca de informações técnicas sobre as especificações de máquinas térmicas e seu funcionamento. Incentivar trabalhos em grupo, com apresentação de resultados.",
    2)

# 3) Replace the "Programa resumido" paragraph text.
$d.Content.Find.Execute(
    "1. Termodinâmica e Energia. 2. Importância das unidades e análise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos domésticos e a Termodinâmica. 5. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 6. Eficiência na conversão de energia. 7. Processos e ciclos térmicos. 8. Termodinâmica e o meio ambiente.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "1. Termodinâmica e Energia. 2. Propriedades das substâncias puras 3. Equipamentos domésticos e a Termodinâmica. 4. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 5. Eficiência na conversão de energia. 6. Processos e ciclos térmicos: equipamentos, materiais e sistemas integrados. 7. Termodinâmica e o meio ambiente",
    2)

# 4) Replace the "Método:" evaluation text.
$d.Content.Find.Execute(
    "Serão realizadas 2 avaliações, com questões abrangendo problemas práticos e conceituais. A 1a. avaliação terá peso 1 e a 2a. avaliação terá peso 2. A nota será a média ponderada das 2 avaliações.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Aulas teóricas expositivas com recursos de mídia variados. Serão realizadas pelo menos duas avaliações escritas abrangendo problemas numéricos e conceituais. Trabalhos em grupo abordando problemas práticos também poderão ser solicitados. Serão envidados esforços para viabilizar viagens didáticas a plantas de geração de potência a fim possibilitar aos alunos o contato com ciclos térmicos reais.",
    2)

# 5) Replace the "Critério:" evaluation text.
$d.Content.Find.Execute(
    "Serão aplicadas duas avaliações escritas (P1, com peso 1 e P2, com peso 2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (P1 + P2)/3.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Somente a nota da última avaliação escrita, aplicada ao final do semestre, terá peso 2. As demais provas escritas ou trabalho em grupo terão peso 1. A nota final será a média ponderada dentre as avaliações aplicadas.",
    2)
